$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object "object[,]" 9,20
$arr[0,0] = "ECs"
$arr[0,1] = "Epha4"
$arr[0,2] = "Efnb1"
$arr[0,3] = "ECs"
$arr[0,4] = 2
$arr[0,5] = 0.6666666666666666
$arr[0,6] = 6.708176333333333
$arr[0,7] = 20.124529
$arr[0,8] = 0.4356329228871633
$arr[0,9] = 0.4356329228871633
$arr[0,10] = 3
$arr[0,11] = 1
$arr[0,12] = 10.31211433333333
$arr[0,13] = 30.936343
$arr[0,14] = 0.633340936097251
$arr[0,15] = 0.633340936097251
$arr[0,16] = 69.17548131749412
$arr[0,17] = 622.579331857447
$arr[0,18] = 0.2759041631761375
$arr[0,19] = 0.2759041631761375
$arr[1,0] = "ECs"
$arr[1,1] = "Epha4"
$arr[1,2] = "Efnb1"
$arr[1,3] = "FAPs"
$arr[1,4] = 2
$arr[1,5] = 0.6666666666666666
$arr[1,6] = 6.708176333333333
$arr[1,7] = 20.124529
$arr[1,8] = 0.4356329228871633
$arr[1,9] = 0.4356329228871633
$arr[1,10] = 3
$arr[1,11] = 1
$arr[1,12] = 4.103438
$arr[1,13] = 12.310314
$arr[1,14] = 0.2520215719230645
$arr[1,15] = 0.2520215719230645
$arr[1,16] = 27.52658567690067
$arr[1,17] = 247.739271092106
$arr[1,18] = 0.109788894007462
$arr[1,19] = 0.109788894007462
$arr[2,0] = "ECs"
$arr[2,1] = "Epha4"
$arr[2,2] = "Efnb1"
$arr[2,3] = "sCs"
$arr[2,4] = 2
$arr[2,5] = 0.6666666666666666
$arr[2,6] = 6.708176333333333
$arr[2,7] = 20.124529
$arr[2,8] = 0.4356329228871633
$arr[2,9] = 0.4356329228871633
$arr[2,10] = 3
$arr[2,11] = 1
$arr[2,12] = 1.866538
$arr[2,13] = 5.599614
$arr[2,14] = 0.1146374919796846
$arr[2,15] = 0.1146374919796846
$arr[2,16] = 12.52106603686733
$arr[2,17] = 112.689594331806
$arr[2,18] = 0.04993986570356374
$arr[2,19] = 0.04993986570356373
$arr[3,0] = "FAPs"
$arr[3,1] = "Epha4"
$arr[3,2] = "Efnb1"
$arr[3,3] = "ECs"
$arr[3,4] = 3
$arr[3,5] = 1
$arr[3,6] = 6.789877333333333
$arr[3,7] = 20.369632
$arr[3,8] = 0.4409386339573907
$arr[3,9] = 0.4409386339573907
$arr[3,10] = 3
$arr[3,11] = 1
$arr[3,12] = 10.31211433333333
$arr[3,13] = 30.936343
$arr[3,14] = 0.633340936097251
$arr[3,15] = 0.633340936097251
$arr[3,16] = 70.01799137064178
$arr[3,17] = 630.161922335776
$arr[3,18] = 0.279264487192017
$arr[3,19] = 0.279264487192017
$arr[4,0] = "FAPs"
$arr[4,1] = "Epha4"
$arr[4,2] = "Efnb1"
$arr[4,3] = "FAPs"
$arr[4,4] = 3
$arr[4,5] = 1
$arr[4,6] = 6.789877333333333
$arr[4,7] = 20.369632
$arr[4,8] = 0.4409386339573907
$arr[4,9] = 0.4409386339573907
$arr[4,10] = 3
$arr[4,11] = 1
$arr[4,12] = 4.103438
$arr[4,13] = 12.310314
$arr[4,14] = 0.2520215719230645
$arr[4,15] = 0.2520215719230645
$arr[4,16] = 27.86184066493866
$arr[4,17] = 250.7565659844479
$arr[4,18] = 0.1111260476515504
$arr[4,19] = 0.1111260476515503
$arr[5,0] = "FAPs"
$arr[5,1] = "Epha4"
$arr[5,2] = "Efnb1"
$arr[5,3] = "sCs"
$arr[5,4] = 3
$arr[5,5] = 1
$arr[5,6] = 6.789877333333333
$arr[5,7] = 20.369632
$arr[5,8] = 0.4409386339573907
$arr[5,9] = 0.4409386339573907
$arr[5,10] = 3
$arr[5,11] = 1
$arr[5,12] = 1.866538
$arr[5,13] = 5.599614
$arr[5,14] = 0.1146374919796846
$arr[5,15] = 0.1146374919796846
$arr[5,16] = 12.67356405800533
$arr[5,17] = 114.062076522048
$arr[5,18] = 0.05054809911382346
$arr[5,19] = 0.05054809911382346
$arr[6,0] = "sCs"
$arr[6,1] = "Epha4"
$arr[6,2] = "Efnb1"
$arr[6,3] = "ECs"
$arr[6,4] = 3
$arr[6,5] = 1
$arr[6,6] = 1.900636333333334
$arr[6,7] = 5.701909000000001
$arr[6,8] = 0.1234284431554459
$arr[6,9] = 0.1234284431554459
$arr[6,10] = 3
$arr[6,11] = 1
$arr[6,12] = 10.31211433333333
$arr[6,13] = 30.936343
$arr[6,14] = 0.633340936097251
$arr[6,15] = 0.633340936097251
$arr[6,16] = 19.59957917542078
$arr[6,17] = 176.396212578787
$arr[6,18] = 0.07817228572909644
$arr[6,19] = 0.07817228572909644
$arr[7,0] = "sCs"
$arr[7,1] = "Epha4"
$arr[7,2] = "Efnb1"
$arr[7,3] = "FAPs"
$arr[7,4] = 3
$arr[7,5] = 1
$arr[7,6] = 1.900636333333334
$arr[7,7] = 5.701909000000001
$arr[7,8] = 0.1234284431554459
$arr[7,9] = 0.1234284431554459
$arr[7,10] = 3
$arr[7,11] = 1
$arr[7,12] = 4.103438
$arr[7,13] = 12.310314
$arr[7,14] = 0.2520215719230645
$arr[7,15] = 0.2520215719230645
$arr[7,16] = 7.799143354380667
$arr[7,17] = 70.192290189426
$arr[7,18] = 0.03110663026405209
$arr[7,19] = 0.03110663026405209
$arr[8,0] = "sCs"
$arr[8,1] = "Epha4"
$arr[8,2] = "Efnb1"
$arr[8,3] = "sCs"
$arr[8,4] = 3
$arr[8,5] = 1
$arr[8,6] = 1.900636333333334
$arr[8,7] = 5.701909000000001
$arr[8,8] = 0.1234284431554459
$arr[8,9] = 0.1234284431554459
$arr[8,10] = 3
$arr[8,11] = 1
$arr[8,12] = 1.866538
$arr[8,13] = 5.599614
$arr[8,14] = 0.1146374919796846
$arr[8,15] = 0.1146374919796846
$arr[8,16] = 3.547609940347334
$arr[8,17] = 31.928489463126
$arr[8,18] = 0.01414952716229739
$arr[8,19] = 0.01414952716229739

$ws.Range("A2:T10").Value = $arr
